$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B7 and B8 values (dni change errors)
$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 2

# Fill in row 10 with the "Cadet A" team data
$ws.Range("A10").Value = "Cadet A"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 20
$ws.Range("E10").Value = 6

# Update the active selection to D11
$ws.Range("D11").Select()
